# Swap values in columns A, B, E, F, G, H, Q, R between rows 6/7 and rows 22/23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

function Swap-Rows($ws, $cols, $row1, $row2) {
    foreach ($col in $cols) {
        $rng1 = $ws.Range("$col$row1")
        $rng2 = $ws.Range("$col$row2")
        $v1 = $rng1.Value2
        $v2 = $rng2.Value2
        $rng1.Value2 = $v2
        $rng2.Value2 = $v1
    }
}

Swap-Rows $ws $cols 6 7
Swap-Rows $ws $cols 22 23
